$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 429.15384
$ws.Range("I2").Value = 398.0909
$ws.Range("K2").Value = 398.0909
$ws.Range("M2").Value = -285.0909
$ws.Range("H86").Value = 5594.3
$ws.Range("I86").Value = 4874.875
$ws.Range("K86").Value = 4874.875
$ws.Range("M86").Value = -3751.875
$ws.Range("H89").Value = 5594.3
$ws.Range("I89").Value = 4874.875
$ws.Range("K89").Value = 24374.375
$ws.Range("M89").Value = -18758.375
$ws.Range("H92").Value = 503
$ws.Range("I92").Value = 503
$ws.Range("K92").Value = 503
$ws.Range("M92").Value = 745
$ws.Range("H123").Value = 140780
$ws.Range("J123").Value = 140780
$ws.Range("L123").Value = 140780
$ws.Range("N123").Value = -150580
$ws.Range("H129").Value = 5831
$ws.Range("I129").Value = 3421
$ws.Range("J129").Value = 12458.5
$ws.Range("K129").Value = 10263
$ws.Range("L129").Value = 37375.5
$ws.Range("M129").Value = -5263
$ws.Range("N129").Value = -47375.5
$ws.Range("H135").Value = 750.8421
$ws.Range("J135").Value = 4000
$ws.Range("L135").Value = 36000
$ws.Range("N135").Value = -41070
$ws.Range("H137").Value = 2248.75
$ws.Range("I137").Value = 2242.7778
$ws.Range("K137").Value = 6728.3334
$ws.Range("M137").Value = -4178.3334
$ws.Range("H138").Value = 2887.707
$ws.Range("I138").Value = 1398.1034
$ws.Range("J138").Value = 3504.8286
$ws.Range("K138").Value = 4194.3102
$ws.Range("L138").Value = 10514.4858
$ws.Range("M138").Value = 945.6898000000001
$ws.Range("N138").Value = -20794.4858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2351.5967
$ws.Range("I32").Value = 2096.5933
$ws.Range("K32").Value = 2096.5933
$ws.Range("M32").Value = -1809.5933
$ws.Range("H61").Value = 4215.7144
$ws.Range("I61").Value = 4094.5833
$ws.Range("K61").Value = 4094.5833
$ws.Range("M61").Value = -3882.5833
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = $null
$ws.Range("H74").Value = 12944.333
$ws.Range("I74").Value = 2083.6667
$ws.Range("K74").Value = 2083.6667
$ws.Range("M74").Value = -1209.6667
$ws.Range("H77").Value = 12944.333
$ws.Range("I77").Value = 2083.6667
$ws.Range("K77").Value = 10418.3335
$ws.Range("M77").Value = -6050.333500000001
$ws.Range("H110").Value = 12187.621
$ws.Range("I110").Value = 14944.263
$ws.Range("K110").Value = 14944.263
$ws.Range("M110").Value = -12899.263
$ws.Range("H122").Value = 2694.4
$ws.Range("I122").Value = 2694.4
$ws.Range("K122").Value = 8083.200000000001
$ws.Range("M122").Value = -5633.200000000001
$ws.Range("H132").Value = 3090.353
$ws.Range("I132").Value = 2809.3845
$ws.Range("J132").Value = 4003.5
$ws.Range("K132").Value = 8428.1535
$ws.Range("L132").Value = 12010.5
$ws.Range("M132").Value = -5898.1535
$ws.Range("N132").Value = -17070.5
$ws.Range("H136").Value = 4215.7144
$ws.Range("I136").Value = 4094.5833
$ws.Range("K136").Value = 12283.7499
$ws.Range("M136").Value = -9733.749899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1227.5
$ws.Range("I86").Value = 1303.3334
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 1303.3334
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -180.3334
$ws.Range("N86").Value = -3246
$ws.Range("H89").Value = 1227.5
$ws.Range("I89").Value = 1303.3334
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 6516.666999999999
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -900.6669999999995
$ws.Range("N89").Value = -16232
$ws.Range("H94").Value = 1710.2727
$ws.Range("I94").Value = 1381.3
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 1381.3
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -930.3
$ws.Range("N94").Value = -5902
$ws.Range("H99").Value = 3893.1875
$ws.Range("I99").Value = 1714.7693
$ws.Range("K99").Value = 1714.7693
$ws.Range("M99").Value = -216.7692999999999
$ws.Range("H106").Value = 59995
$ws.Range("J106").Value = 59995
$ws.Range("L106").Value = 59995
$ws.Range("N106").Value = -62519
$ws.Range("H107").Value = 2459.8
$ws.Range("I107").Value = 2459.8
$ws.Range("K107").Value = 2459.8
$ws.Range("M107").Value = -539.8000000000002
$ws.Range("H134").Value = 2874.389
$ws.Range("I134").Value = 2590.2856
$ws.Range("K134").Value = 7770.8568
$ws.Range("M134").Value = -5235.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22442.766
$ws.Range("I31").Value = 24423.137
$ws.Range("K31").Value = 24423.137
$ws.Range("M31").Value = -24128.137
$ws.Range("H34").Value = 22442.766
$ws.Range("I34").Value = 24423.137
$ws.Range("K34").Value = 24423.137
$ws.Range("M34").Value = -24221.137
$ws.Range("H58").Value = 3014.1667
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = $null
$ws.Range("H107").Value = 592.94446
$ws.Range("I107").Value = 601.13794
$ws.Range("J107").Value = 559
$ws.Range("K107").Value = 601.13794
$ws.Range("L107").Value = 559
$ws.Range("M107").Value = 1318.86206
$ws.Range("N107").Value = -4399
$ws.Range("H136").Value = 3014.1667
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null
$ws.Range("H141").Value = 64371.43
$ws.Range("J141").Value = 64371.43
$ws.Range("L141").Value = 64371.43
$ws.Range("N141").Value = -74731.42999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1368.8966
$ws.Range("I132").Value = 1133.25
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 10199.25
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -7669.25
$ws.Range("N132").Value = -27560
$ws.Range("H134").Value = 5323.421
$ws.Range("I134").Value = 3134.0625
$ws.Range("K134").Value = 9402.1875
$ws.Range("M134").Value = -4332.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15255.272
$ws.Range("I70").Value = 14597.625
$ws.Range("K70").Value = 14597.625
$ws.Range("M70").Value = -14327.625
$ws.Range("H73").Value = 15255.272
$ws.Range("I73").Value = 14597.625
$ws.Range("K73").Value = 14597.625
$ws.Range("M73").Value = -13661.625
$ws.Range("H97").Value = 776.7619
$ws.Range("I97").Value = 716.75
$ws.Range("J97").Value = 856.7778
$ws.Range("K97").Value = 716.75
$ws.Range("L97").Value = 856.7778
$ws.Range("M97").Value = -220.75
$ws.Range("N97").Value = -1848.7778
$ws.Range("H107").Value = 1423.625
$ws.Range("I107").Value = 1177.5555
$ws.Range("J107").Value = 1740
$ws.Range("K107").Value = 1177.5555
$ws.Range("L107").Value = 1740
$ws.Range("M107").Value = 742.4445000000001
$ws.Range("N107").Value = -5580
$ws.Range("H122").Value = 2275.0833
$ws.Range("I122").Value = 2254.7778
$ws.Range("J122").Value = 2336
$ws.Range("K122").Value = 6764.3334
$ws.Range("L122").Value = 7008
$ws.Range("M122").Value = -4314.3334
$ws.Range("N122").Value = -11908
$ws.Range("H132").Value = 191443.25
$ws.Range("I132").Value = 202855.66
$ws.Range("K132").Value = 608566.98
$ws.Range("M132").Value = -606036.98

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = $null
$ws.Range("H105").Value = 29975
$ws.Range("J105").Value = 29975
$ws.Range("L105").Value = 29975
$ws.Range("N105").Value = -36963
